$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19, shifting existing rows 19-21 down to 20-22
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the new record
$ws.Cells.Item(19, 1).Value = 4
$ws.Cells.Item(19, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(19, 3).Value = "Los Lagos"
$ws.Cells.Item(19, 4).Value = 44722
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 10
$ws.Cells.Item(19, 6).Value = 100112012
$ws.Cells.Item(19, 7).Value = "Espinaca"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 30
$ws.Cells.Item(19, 11).Value = 13000
$ws.Cells.Item(19, 12).Value = 13000
$ws.Cells.Item(19, 13).Value = 13000
$ws.Cells.Item(19, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 1300
$ws.Cells.Item(19, 17).Value = 10
$ws.Cells.Item(19, 18).Value = "Hortaliza"
